$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must stay text (force Text format
# first, matching how the source data was authored as strings).
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D18", "D19", "D20", "D21", "D22", "D23", "D27", "D29", "D33", "D34", "D39", "D40", "D41", "D42", "D43", "D44", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '67.844.73'
$ws.Cells.Item(2, 5).Value = '  +1.21%  '
$ws.Cells.Item(3, 4).Value = '3.506.42'
$ws.Cells.Item(3, 5).Value = '  -0.19%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).Value = '599.95'
$ws.Cells.Item(5, 5).Value = '  +0.91%  '
$ws.Cells.Item(6, 4).Value = '182.72'
$ws.Cells.Item(6, 5).Value = '  +5.34%  '
$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$ws.Cells.Item(8, 4).Value = '0.599'
$ws.Cells.Item(8, 5).Value = '  +1.11%  '
$ws.Cells.Item(9, 4).Value = '0.139'
$ws.Cells.Item(9, 5).Value = '  +4.09%  '
$ws.Cells.Item(10, 4).Value = '7.13'
$ws.Cells.Item(10, 5).Value = '  -1.90%  '
$ws.Cells.Item(11, 4).Value = '0.435'
$ws.Cells.Item(11, 5).Value = '  -0.10%  '
$ws.Cells.Item(12, 4).Value = '4.112.95'
$ws.Cells.Item(12, 5).Value = '  -0.17%  '
$ws.Cells.Item(13, 4).Value = '32.30'
$ws.Cells.Item(13, 5).Value = '  +12.49%  '
$ws.Cells.Item(14, 5).Value = '  -0.04%  '
$ws.Cells.Item(15, 4).Value = '67.815.85'
$ws.Cells.Item(15, 5).Value = '  +1.20%  '
$ws.Cells.Item(16, 5).Value = '  +0.77%  '
$ws.Cells.Item(17, 4).Value = '3.506.07'
$ws.Cells.Item(17, 5).Value = '  +0.30%  '
$ws.Cells.Item(18, 4).Value = '6.39'
$ws.Cells.Item(18, 5).Value = '  +1.12%  '
$ws.Cells.Item(19, 4).Value = '14.76'
$ws.Cells.Item(19, 5).Value = '  +3.81%  '
$ws.Cells.Item(20, 4).Value = '396.77'
$ws.Cells.Item(20, 5).Value = '  +0.49%  '
$ws.Cells.Item(21, 4).Value = '8.09'
$ws.Cells.Item(21, 5).Value = '  +1.60%  '
$ws.Cells.Item(22, 4).Value = '73.47'
$ws.Cells.Item(22, 5).Value = '  +0.27%  '
$ws.Cells.Item(23, 4).Value = '0.546'
$ws.Cells.Item(23, 5).Value = '  +1.16%  '
$ws.Cells.Item(24, 5).Value = '  +0.12%  '
$ws.Cells.Item(25, 5).Value = '  +0.10%  '
$ws.Cells.Item(26, 5).Value = '  +2.52%  '
$ws.Cells.Item(27, 4).Value = '10.60'
$ws.Cells.Item(27, 5).Value = '  +4.39%  '
$ws.Cells.Item(28, 5).Value = '  -0.73%  '
$ws.Cells.Item(29, 4).Value = '0.996'
$ws.Cells.Item(29, 5).Value = '  -0.26%  '
$ws.Cells.Item(30, 5).Value = '  +0.30%  '
$ws.Cells.Item(31, 5).Value = '  +1.29%  '
$ws.Cells.Item(32, 5).Value = '  +0.22%  '
$ws.Cells.Item(33, 4).Value = '24.15'
$ws.Cells.Item(33, 5).Value = '  +0.84%  '
$ws.Cells.Item(34, 4).Value = '7.46'
$ws.Cells.Item(34, 5).Value = '  +1.21%  '
$ws.Cells.Item(35, 5).Value = '  +0.09%  '
$ws.Cells.Item(36, 5).Value = '  +2.07%  '
$ws.Cells.Item(37, 5).Value = '  +0.77%  '
$ws.Cells.Item(38, 5).Value = '  +3.20%  '
$ws.Cells.Item(39, 4).Value = '0.876'
$ws.Cells.Item(39, 5).Value = '  -2.21%  '
$ws.Cells.Item(40, 4).Value = '7.16'
$ws.Cells.Item(40, 5).Value = '  +3.80%  '
$ws.Cells.Item(41, 4).Value = '4.76'
$ws.Cells.Item(41, 5).Value = '  +1.82%  '
$ws.Cells.Item(42, 4).Value = '2.72'
$ws.Cells.Item(42, 5).Value = '  +4.09%  '
$ws.Cells.Item(43, 4).Value = '27.86'
$ws.Cells.Item(43, 5).Value = '  +2.92%  '
$ws.Cells.Item(44, 4).Value = '26.75'
$ws.Cells.Item(44, 5).Value = '  +1.43%  '
$ws.Cells.Item(45, 5).Value = '  -0.42%  '
$ws.Cells.Item(46, 4).Value = '2.842.08'
$ws.Cells.Item(46, 5).Value = '  +1.56%  '
$ws.Cells.Item(47, 5).Value = '  -1.06%  '
$ws.Cells.Item(48, 5).Value = '  +0.49%  '
$ws.Cells.Item(49, 4).Value = '345.11'
$ws.Cells.Item(49, 5).Value = '  +1.30%  '
$ws.Cells.Item(50, 4).Value = '1.09'
$ws.Cells.Item(50, 5).Value = '  -0.27%  '
$ws.Cells.Item(51, 4).Value = '33.75'
$ws.Cells.Item(51, 5).Value = '  +0.86%  '
